# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment schedule"
# sheet to make room for a Variable Instalments column, then switch the
# active sheet/selection to "Repayment schedule".

$wb = $excel.ActiveWorkbook

$repaymentSheet = $wb.Worksheets.Item("Repayment schedule")
$repaymentSheet.Columns("N:N").Insert()
$repaymentSheet.Columns("N:N").ColumnWidth = $repaymentSheet.Columns("M:M").ColumnWidth

$repaymentSheet.Activate()
$repaymentSheet.Range("S8").Select()
